# Add a new slide (6th) using the same "Title and Content" layout as the
# other content slides in the deck (slideLayout2.xml == ppLayoutText).
# The new slide is left blank (title + content placeholders with no text),
# matching the freshly-inserted slide added to the deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Add($p.Slides.Count + 1, 2)
